$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $ref, $proto, $text)
    $ws.Range($proto).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
    $ws.Range($ref).Value = $text
}

# Row 10
Set-CellText $ws 'A10' 'A3' 'Objetivos:'
Set-CellText $ws 'B10' 'B3' 'Promover a formação do engenheiro de materiais sob o ponto de vista do desenvolvimento de competências gerais e específicas.Integrar essa disciplina com outras de semestres anteriores e do presente semestre da matriz curricular do curso de Engenharia de Materiais. No início do semestre o responsável pela disciplina deve reunir com os professores de semestres anteriores e presente para planejar trabalhos conjuntos. O tema do trabalho de projeto será definido em conjunto os professores responsáveis pelas disciplinas envolvidas, onde um dos professores será o orientador do respectivo projeto. - Incentivar trabalhos em grupo, com apresentação de resultados.'
Set-CellText $ws 'C10' 'C3' 'Promover a formação do engenheiro de materiais sob o ponto de vista do desenvolvimento de competências gerais e específicas.Integrar essa disciplina com outras de semestres anteriores e do presente semestre da matriz curricular do curso de Engenharia de Materiais. No início do semestre o responsável pela disciplina deve reunir com os professores de semestres anteriores e presente para planejar trabalhos conjuntos. O tema do trabalho de projeto será definido em conjunto os professores responsáveis pelas disciplinas envolvidas, onde um dos professores será o orientador do respectivo projeto. - Incentivar trabalhos em grupo, com apresentação de resultados.'

# Row 11
Set-CellText $ws 'A11' 'A3' 'Objectives:'

# Row 12
Set-CellText $ws 'A12' 'A3' 'Docentes responsáveis:'

# Row 13
Set-CellText $ws 'B13' 'B3' '7459752 - Maria Ismenia Sodero Toledo Faria'
Set-CellText $ws 'C13' 'C3' '7459752 - Maria Ismenia Sodero Toledo Faria'

# Row 14
Set-CellText $ws 'B14' 'B3' '2166002 - Sandra Giacomin Schneider'
Set-CellText $ws 'C14' 'C3' '2166002 - Sandra Giacomin Schneider'

# Row 15
Set-CellText $ws 'B15' 'B3' '1922320 - Sebastiao Ribeiro'
Set-CellText $ws 'C15' 'C3' '1922320 - Sebastiao Ribeiro'

# Row 16
Set-CellText $ws 'A16' 'A3' 'Programa resumido:'
Set-CellText $ws 'B16' 'B3' 'Entender as principais abordagens para o desenvolvimento de produtos. Definir, planejar e projetar modelo de negócios inovadores. Conhecer e aplicar os principais modelos de gestão e operação para startups. Aplicar técnicas para modelagem financeira de novos empreendimentos. Aplicar conceitos de marketing para novos negócios.'
Set-CellText $ws 'C16' 'C3' 'Entender as principais abordagens para o desenvolvimento de produtos. Definir, planejar e projetar modelo de negócios inovadores. Conhecer e aplicar os principais modelos de gestão e operação para startups. Aplicar técnicas para modelagem financeira de novos empreendimentos. Aplicar conceitos de marketing para novos negócios.'

# Row 17
Set-CellText $ws 'A17' 'A3' 'Short syllabus:'

# Row 18
Set-CellText $ws 'A18' 'A3' 'Programa:'
Set-CellText $ws 'B18' 'B3' '1. Modelo de negócios: conceitos, cases, abordagens de projeto de modelos de negócios, operações em Marketing digital, modelos de financiamento de startups, tipos de investidores, valoração do empreendimento e decisões de saída2. Planejamento e pesquisa sobre tecnologias visando inovação tecnológica3. Planejamento e pesquisa sobre modelos de negócio, tipologias e arquiteturas 4. Visitas a incubadoras e aceleradoras de startups 5. Desenvolvimento e apresentação de um pitch de negócio6. Elaboração de um plano de negócio'
Set-CellText $ws 'C18' 'C3' '1. Modelo de negócios: conceitos, cases, abordagens de projeto de modelos de negócios, operações em Marketing digital, modelos de financiamento de startups, tipos de investidores, valoração do empreendimento e decisões de saída2. Planejamento e pesquisa sobre tecnologias visando inovação tecnológica3. Planejamento e pesquisa sobre modelos de negócio, tipologias e arquiteturas 4. Visitas a incubadoras e aceleradoras de startups 5. Desenvolvimento e apresentação de um pitch de negócio6. Elaboração de um plano de negócio'

# Row 19
Set-CellText $ws 'A19' 'A3' 'Syllabus:'

# Row 20
Set-CellText $ws 'A20' 'A3' 'Avaliação:'

# Row 21
Set-CellText $ws 'A21' 'A3' 'Método:'
Set-CellText $ws 'B21' 'B3' 'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.'
Set-CellText $ws 'C21' 'C3' 'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.'

# Row 22
Set-CellText $ws 'A22' 'A3' 'Critério:'
Set-CellText $ws 'B22' 'B3' 'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.'
Set-CellText $ws 'C22' 'C3' 'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.'

# Row 23
Set-CellText $ws 'A23' 'A3' 'Norma de recuperação:'
Set-CellText $ws 'B23' 'B3' 'não há'
Set-CellText $ws 'C23' 'C3' 'não há'

# Row 24
Set-CellText $ws 'A24' 'A3' 'Bibliografia:'
Set-CellText $ws 'B24' 'B3' '- MCCAHAN, S.; ANDERSON, P.; KORTSCHOT, M.; WEISS, P.; WOODHOUSE, K. Projetos de Engenharia: uma introdução. 1ª edição. -Rio de Janeiro: LTC, 2017.- BROCKMAN, Jay B. Introdução à Engenharia - Modelagem e solução de problemas. Rio de Janeiro: LTC, 2010.- CAVALCANTI, Carolina C.; FILATRO, Andrea C. Design Thinking na educação presencial, a distância e corporativa. São Paulo: Editora Saraiva, 2016.- FINOCCHIO, José. PMC Projeto modelo Canvas, 3 ed. São Paulo: Editora Saraiva, 2020.- CAMARGO, Robson; RIBAS, Thomaz. Gestão ágil de projetos: As melhores soluções para suas necessidades. São Paulo: Editora Saraiva, 2019.- BRANCO, R. H. F.; LEITE, D. E.; VINHA JR., Rubem. Gestão colaborativa de projetos: A combinação de Design Thinking e ferramentas práticas para gerenciar seus projetos. São Paulo: Editora Saraiva Universitária, 2016- OSTERWALDER, Alexander; PIGNEUR, Yves. Business Model Generation: Inovação em modelos de negócios. Rio de Janeiro: Alta Books, 2011.'
Set-CellText $ws 'C24' 'C3' '- MCCAHAN, S.; ANDERSON, P.; KORTSCHOT, M.; WEISS, P.; WOODHOUSE, K. Projetos de Engenharia: uma introdução. 1ª edição. -Rio de Janeiro: LTC, 2017.- BROCKMAN, Jay B. Introdução à Engenharia - Modelagem e solução de problemas. Rio de Janeiro: LTC, 2010.- CAVALCANTI, Carolina C.; FILATRO, Andrea C. Design Thinking na educação presencial, a distância e corporativa. São Paulo: Editora Saraiva, 2016.- FINOCCHIO, José. PMC Projeto modelo Canvas, 3 ed. São Paulo: Editora Saraiva, 2020.- CAMARGO, Robson; RIBAS, Thomaz. Gestão ágil de projetos: As melhores soluções para suas necessidades. São Paulo: Editora Saraiva, 2019.- BRANCO, R. H. F.; LEITE, D. E.; VINHA JR., Rubem. Gestão colaborativa de projetos: A combinação de Design Thinking e ferramentas práticas para gerenciar seus projetos. São Paulo: Editora Saraiva Universitária, 2016- OSTERWALDER, Alexander; PIGNEUR, Yves. Business Model Generation: Inovação em modelos de negócios. Rio de Janeiro: Alta Books, 2011.'

# Row 25
Set-CellText $ws 'A25' 'A3' 'Requisitos:'

# Row 26
Set-CellText $ws 'B26' 'B3' 'LOM3108 -  Projeto Integrado em Engenharia de Materiais II  (Requisito fraco)
'
Set-CellText $ws 'C26' 'C3' 'LOM3108 -  Projeto Integrado em Engenharia de Materiais II  (Requisito fraco)
'

# Clear stale cells left over from the previous layout
$ws.Range('A13').Clear() | Out-Null
$ws.Range('A14').Clear() | Out-Null
$ws.Range('A15').Clear() | Out-Null
$ws.Range('B19').Clear() | Out-Null
$ws.Range('C19').Clear() | Out-Null
$ws.Range('B20').Clear() | Out-Null
$ws.Range('C20').Clear() | Out-Null

# Row heights: rows that lose their custom height go back to auto
$ws.Rows.Item(11).AutoFit() | Out-Null
$ws.Rows.Item(13).AutoFit() | Out-Null
$ws.Rows.Item(14).AutoFit() | Out-Null
$ws.Rows.Item(15).AutoFit() | Out-Null
$ws.Rows.Item(20).AutoFit() | Out-Null

# Row heights: rows with an explicit custom height in the final layout
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).RowHeight = 120
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 60
$ws.Rows.Item(24).RowHeight = 120
$ws.Rows.Item(26).RowHeight = 30

